# Update raw data dates in EventData sheet (02 Nov 2011) and sheet view/selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EventData")

# Column I, rows 2-51: change date serial 40834 (18-Oct-2011) to 40849 (02-Nov-2011)
for ($r = 2; $r -le 51; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    if ($cell.Value2 -eq 40834) {
        $cell.Value2 = 40849
    }
}

# Update the sheet view's top-left cell and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 50
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K5").Select()
